# Generate Report for Handoff
#
# The "b.md" row (row 3) on the Overview sheet, and on each language
# sheet (zh-cn, de-de), is refreshed to reflect that a new handoff
# package was generated for b.md:
#   - Status moves from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   - Content Duplicate flips from True to False
#   - A new Latest Handoff File (the "b.*.xlf") and a new Latest
#     Handoff Datetime are recorded
#   - An Error Detail message is attached, noting the handback file is
#     stale relative to the newly generated handoff
#   - The Error Detail column is widened so the message is readable

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet - row 3 is the b.md entry
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 06:41:08"

# ---------------------------------------------------------------
# zh-cn sheet - row 3 is the b.md entry
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps this a text "False" (matches the source
# data type) instead of letting it be auto-detected as a Boolean.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-04 06:40:58"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fe956781b323be44c1014d25f60498408a84bba/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/626ba5d153c043bc9f988d9c497208086debf95f/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.2

# ---------------------------------------------------------------
# de-de sheet - row 3 is the b.md entry
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-04 06:41:08"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fe956781b323be44c1014d25f60498408a84bba/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/626ba5d153c043bc9f988d9c497208086debf95f/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.2
